# Natmi following Dr Hou advice
# Expand LR-pair result table from 4 rows (one per sending cluster) to 8 rows
# (one per sending cluster x target cluster in {M2, sCs}).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for 4 additional data rows (rows 6-9)
$ws.Range("A6:A9").EntireRow.Insert()

# Ligand/Receptor symbol columns are constant for every row
$ws.Range("B2:B9").Value = "Col3a1"
$ws.Range("C2:C9").Value = "Mag"

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "M2"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.822558333333333
$ws.Range("H2").Value = 23.467675
$ws.Range("I2").Value = 0.001247993910151231
$ws.Range("J2").Value = 0.001247993910151231
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3133113333333333
$ws.Range("N2").Value = 0.939934
$ws.Range("O2").Value = 0.4010297802586483
$ws.Range("P2").Value = 0.4010297802586483
$ws.Range("Q2").Value = 2.450896181494445
$ws.Range("R2").Value = 22.05806563345
$ws.Range("S2").Value = 0.0005004827235520794
$ws.Range("T2").Value = 0.0005004827235520794

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.822558333333333
$ws.Range("H3").Value = 23.467675
$ws.Range("I3").Value = 0.001247993910151231
$ws.Range("J3").Value = 0.001247993910151231
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4679556666666667
$ws.Range("N3").Value = 1.403867
$ws.Range("O3").Value = 0.5989702197413518
$ws.Range("P3").Value = 0.5989702197413518
$ws.Range("Q3").Value = 3.660610499913889
$ws.Range("R3").Value = 32.945494499225
$ws.Range("S3").Value = 0.0007475111865991517
$ws.Range("T3").Value = 0.0007475111865991517

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6037.773437333333
$ws.Range("H4").Value = 18113.320312
$ws.Range("I4").Value = 0.9632532171165058
$ws.Range("J4").Value = 0.9632532171165058
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3133113333333333
$ws.Range("N4").Value = 0.939934
$ws.Range("O4").Value = 0.4010297802586483
$ws.Range("P4").Value = 0.4010297802586483
$ws.Range("Q4").Value = 1891.70284601549
$ws.Range("R4").Value = 17025.32561413941
$ws.Range("S4").Value = 0.3862932259936683
$ws.Range("T4").Value = 0.3862932259936683

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6037.773437333333
$ws.Range("H5").Value = 18113.320312
$ws.Range("I5").Value = 0.9632532171165058
$ws.Range("J5").Value = 0.9632532171165058
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.4679556666666667
$ws.Range("N5").Value = 1.403867
$ws.Range("O5").Value = 0.5989702197413518
$ws.Range("P5").Value = 0.5989702197413518
$ws.Range("Q5").Value = 2825.410294049611
$ws.Range("R5").Value = 25428.6926464465
$ws.Range("S5").Value = 0.5769599911228375
$ws.Range("T5").Value = 0.5769599911228375

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.018774
$ws.Range("H6").Value = 3.056322
$ws.Range("I6").Value = 0.0001625329839219791
$ws.Range("J6").Value = 0.0001625329839219791
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3133113333333333
$ws.Range("N6").Value = 0.939934
$ws.Range("O6").Value = 0.4010297802586483
$ws.Range("P6").Value = 0.4010297802586483
$ws.Range("Q6").Value = 0.3191934403053333
$ws.Range("R6").Value = 2.872740962748
$ws.Range("S6").Value = 0.00006518056682701368
$ws.Range("T6").Value = 0.00006518056682701368

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.018774
$ws.Range("H7").Value = 3.056322
$ws.Range("I7").Value = 0.0001625329839219791
$ws.Range("J7").Value = 0.0001625329839219791
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.4679556666666667
$ws.Range("N7").Value = 1.403867
$ws.Range("O7").Value = 0.5989702197413518
$ws.Range("P7").Value = 0.5989702197413518
$ws.Range("Q7").Value = 0.4767410663526666
$ws.Range("R7").Value = 4.290669597173999
$ws.Range("S7").Value = 0.0000973524170949654
$ws.Range("T7").Value = 0.0000973524170949654

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 221.4914043333333
$ws.Range("H8").Value = 664.474213
$ws.Range("I8").Value = 0.03533625598942085
$ws.Range("J8").Value = 0.03533625598942085
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3133113333333333
$ws.Range("N8").Value = 0.939934
$ws.Range("O8").Value = 0.4010297802586483
$ws.Range("P8").Value = 0.4010297802586483
$ws.Range("Q8").Value = 69.3957672135491
$ws.Range("R8").Value = 624.561904921942
$ws.Range("S8").Value = 0.01417089097460079
$ws.Range("T8").Value = 0.01417089097460079

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 221.4914043333333
$ws.Range("H9").Value = 664.474213
$ws.Range("I9").Value = 0.03533625598942085
$ws.Range("J9").Value = 0.03533625598942085
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4679556666666667
$ws.Range("N9").Value = 1.403867
$ws.Range("O9").Value = 0.5989702197413518
$ws.Range("P9").Value = 0.5989702197413518
$ws.Range("Q9").Value = 103.6481577757412
$ws.Range("R9").Value = 932.8334199816709
$ws.Range("S9").Value = 0.02116536501482006
$ws.Range("T9").Value = 0.02116536501482006
